$wb = $excel.ActiveWorkbook

# Helper: write a digit-only string into a cell while keeping it a true TEXT
# value (not auto-coerced to a number) and WITHOUT leaving a "Text" number
# format behind on the cell. We stage the text in a scratch cell (formatted
# as Text so Excel won't re-parse it as a number), then copy/paste only the
# VALUE (xlPasteValues) onto the real target, which carries the text over
# without carrying the scratch cell's number format along with it.
function Set-TextValue {
    param($targetRange, [string]$text)
    $sheet = $targetRange.Worksheet
    $scratch = $sheet.Range("ZZ9999")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $targetRange.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# Existing sheet is "ODI Batting" (sheet1). We need to insert a new sheet
# "Player Info" BEFORE it, becoming the first sheet.
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

# The act of inserting a sheet shifts identities around in this runtime,
# so re-resolve "ODI Batting" by name before touching it again.
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# Fill Player Info sheet
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Copy the bold/bordered header formatting used by the other sheet's header row
$odiBatting.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

Set-TextValue $playerInfo.Range("A2") "6607"
$playerInfo.Range("B2").Value = "Tristan Stubbs"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# Update ODI Batting sheet: rename MATCH_CARD_LINK -> MATCH_CODE, and change the URL value to just the match code
$odiBatting.Range("D1").Value = "MATCH_CODE"
Set-TextValue $odiBatting.Range("D2") "4727"
